$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 212.5
$ws.Range("I39").Value = 9.428572000000001
$ws.Range("J39").Value = 686.3333
$ws.Range("K39").Value = 28.285716
$ws.Range("L39").Value = 2058.9999
$ws.Range("M39").Value = 267.714284
$ws.Range("N39").Value = -2650.9999
$ws.Range("H74").Value = 7489.05
$ws.Range("I74").Value = 7646.3687
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 7646.3687
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -6710.3687
$ws.Range("N74").Value = -6372
$ws.Range("H76").Value = 4512.2856
$ws.Range("I76").Value = 4631
$ws.Range("J76").Value = 3800
$ws.Range("K76").Value = 4631
$ws.Range("L76").Value = 3800
$ws.Range("M76").Value = -4316
$ws.Range("N76").Value = -4430
$ws.Range("H77").Value = 7489.05
$ws.Range("I77").Value = 7646.3687
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 38231.8435
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -33551.8435
$ws.Range("N77").Value = -31860
$ws.Range("H79").Value = 4512.2856
$ws.Range("I79").Value = 4631
$ws.Range("J79").Value = 3800
$ws.Range("K79").Value = 4631
$ws.Range("L79").Value = 3800
$ws.Range("M79").Value = -3539
$ws.Range("N79").Value = -5984
$ws.Range("H80").Value = 2106.9092
$ws.Range("I80").Value = 678.8
$ws.Range("J80").Value = 3297
$ws.Range("K80").Value = 2036.4
$ws.Range("L80").Value = 9891
$ws.Range("M80").Value = -1038.4
$ws.Range("N80").Value = -11887
$ws.Range("H83").Value = 2106.9092
$ws.Range("I83").Value = 678.8
$ws.Range("J83").Value = 3297
$ws.Range("K83").Value = 6109.2
$ws.Range("L83").Value = 29673
$ws.Range("M83").Value = -1117.2
$ws.Range("N83").Value = -39657
$ws.Range("H113").Value = 9124.25
$ws.Range("I113").Value = 8499
$ws.Range("H137").Value = 14494107
$ws.Range("I137").Value = 20834372
$ws.Range("J137").Value = 2071.1428
$ws.Range("K137").Value = 62503116
$ws.Range("L137").Value = 6213.428400000001
$ws.Range("M137").Value = -62500566
$ws.Range("N137").Value = -11313.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1778.3334
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2335
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2335
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3147
$ws.Range("H91").Value = 1778.3334
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2335
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2335
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5143
$ws.Range("H122").Value = 4351.6113
$ws.Range("I122").Value = 4351.6113
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13054.8339
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10604.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5227.25
$ws.Range("I105").Value = 5227.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5227.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3480.25
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("H134").Value = 1826
$ws.Range("I134").Value = 1826
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5478
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2943
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62500210
$ws.Range("I7").Value = 66666890
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 66666890
$ws.Range("L7").Value = 35
$ws.Range("M7").Value = -66666777
$ws.Range("N7").Value = -261
$ws.Range("H31").Value = 2671.4211
$ws.Range("I31").Value = 2266.5
$ws.Range("J31").Value = 2965.9092
$ws.Range("K31").Value = 2266.5
$ws.Range("L31").Value = 2965.9092
$ws.Range("M31").Value = -1971.5
$ws.Range("N31").Value = -3555.9092
$ws.Range("H34").Value = 2671.4211
$ws.Range("I34").Value = 2266.5
$ws.Range("J34").Value = 2965.9092
$ws.Range("K34").Value = 2266.5
$ws.Range("L34").Value = 2965.9092
$ws.Range("M34").Value = -2064.5
$ws.Range("N34").Value = -3369.9092
$ws.Range("H62").Value = 2179.8
$ws.Range("I62").Value = 1966.6666
$ws.Range("J62").Value = 2499.5
$ws.Range("K62").Value = 1966.6666
$ws.Range("L62").Value = 2499.5
$ws.Range("M62").Value = -1342.6666
$ws.Range("N62").Value = -3747.5
$ws.Range("H65").Value = 2179.8
$ws.Range("I65").Value = 1966.6666
$ws.Range("J65").Value = 2499.5
$ws.Range("K65").Value = 9833.333000000001
$ws.Range("L65").Value = 12497.5
$ws.Range("M65").Value = -6713.333000000001
$ws.Range("N65").Value = -18737.5
$ws.Range("H99").Value = 7700
$ws.Range("I99").Value = 5323.1665
$ws.Range("J99").Value = 10552.2
$ws.Range("K99").Value = 5323.1665
$ws.Range("L99").Value = 10552.2
$ws.Range("M99").Value = -3825.1665
$ws.Range("N99").Value = -13548.2
$ws.Range("H126").Value = 7700
$ws.Range("I126").Value = 5323.1665
$ws.Range("J126").Value = 10552.2
$ws.Range("K126").Value = 15969.4995
$ws.Range("L126").Value = 31656.6
$ws.Range("M126").Value = -13499.4995
$ws.Range("N126").Value = -36596.60000000001
$ws.Range("H134").Value = 2541.1333
$ws.Range("I134").Value = 2547.4614
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 7642.3842
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -5107.3842
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3099061.2
$ws.Range("I4").Value = 857316.25
$ws.Range("J4").Value = 30000000
$ws.Range("K4").Value = 2571948.75
$ws.Range("L4").Value = 90000000
$ws.Range("M4").Value = -2571836.75
$ws.Range("N4").Value = -90000224
$ws.Range("H47").Value = 3554.875
$ws.Range("I47").Value = 3554.875
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 10664.625
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -10233.625
$ws.Range("H113").Value = 2679.889
$ws.Range("I113").Value = 2563.2
$ws.Range("J113").Value = 2825.75
$ws.Range("K113").Value = 7689.599999999999
$ws.Range("L113").Value = 8477.25
$ws.Range("M113").Value = -5519.599999999999
$ws.Range("N113").Value = -12817.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998.5
$ws.Range("I70").Value = 4998.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4998.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4728.5
$ws.Range("H73").Value = 4998.5
$ws.Range("I73").Value = 4998.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4998.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4062.5
$ws.Range("H80").Value = 2620.6667
$ws.Range("I80").Value = 2440.8572
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 2440.8572
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -1442.8572
$ws.Range("N80").Value = -5246
$ws.Range("H83").Value = 2620.6667
$ws.Range("I83").Value = 2440.8572
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 12204.286
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -7212.286
$ws.Range("N83").Value = -26234
$ws.Range("H102").Value = 1899.0358
$ws.Range("I102").Value = 1891.2693
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1891.2693
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -269.2692999999999
$ws.Range("N102").Value = -5244
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 50000
$ws.Range("I94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("K94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("M94").Value = -49099
$ws.Range("N94").Value = -51802
